# updated inventory and rookie decks
$wb = $excel.ActiveWorkbook

$blue = $wb.Worksheets.Item("blue")
$red  = $wb.Worksheets.Item("red")

# --- blue sheet: Disperse (row 3) and Seacoast Drake (row 10) now fully bought in ---
$blue.Range("C3").Value = 2
$blue.Range("C10").Value = 4

# --- red sheet: Goblin Arsonist (row 8) and Thundering Giant (row 11) now fully bought in ---
$red.Range("C8").Value = 4
$red.Range("C11").Value = 4

# red's card_name column had no explicit width before; size it to fit the
# longest card name, same as the blue sheet already does
$red.Columns.Item(1).ColumnWidth = 18.85546875

# keep blue as the active tab/selection, now resting on C11
$blue.Activate()
$blue.Range("C11").Select()
